# Revert ""Fuel" sheet update for both CH and SIN"
# This reverts the prior commit's changes to the FUELS sheet of the
# LCA_infrastructure.xlsx workbook (SIN database).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FUELS")

# C2: was a formula (=1.1767+0.0019487+0.0000015726) -> plain literal value
$ws.Range("C2").Value = 1.403

# D2: literal value change
$ws.Range("D2").Value = 0.1

# F2: reference text reverts to the same "from CEA, costs in USD-2015," note
# already used by F3/F4 (drops the bespoke ecoinvent natural-gas string)
$ws.Range("F2").Value = "from CEA, costs in USD-2015,"

# Restore the previously-selected cell on this sheet
$ws.Range("F16").Select()
